$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 86815760
$ws.Range("I116").Value = 135435410
$ws.Range("J116").Value = 62505932
$ws.Range("K116").Value = 135435410
$ws.Range("L116").Value = 62505932
$ws.Range("M116").Value = -135431968
$ws.Range("N116").Value = -62512816
$ws.Range("H131").Value = 1624.2142
$ws.Range("I131").Value = 1145.0834
$ws.Range("J131").Value = 4499
$ws.Range("K131").Value = 3435.2502
$ws.Range("L131").Value = 13497
$ws.Range("M131").Value = 1604.7498
$ws.Range("N131").Value = -23577
$ws.Range("H137").Value = 3545
$ws.Range("I137").Value = 1212.7894
$ws.Range("J137").Value = 6006.778
$ws.Range("K137").Value = 3638.3682
$ws.Range("L137").Value = 18020.334
$ws.Range("M137").Value = -1088.3682
$ws.Range("N137").Value = -23120.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4954.654
$ws.Range("I32").Value = 3089.6316
$ws.Range("K32").Value = 3089.6316
$ws.Range("M32").Value = -2802.6316
$ws.Range("H45").Value = 1999.5
$ws.Range("I45").Value = 1999.5
$ws.Range("K45").Value = 1999.5
$ws.Range("M45").Value = -1622.5
$ws.Range("H61").Value = 13434.5
$ws.Range("I61").Value = 2799
$ws.Range("K61").Value = 2799
$ws.Range("M61").Value = -2587
$ws.Range("H122").Value = 4639.057
$ws.Range("I122").Value = 2082.3635
$ws.Range("J122").Value = 8965.77
$ws.Range("K122").Value = 6247.0905
$ws.Range("L122").Value = 26897.31
$ws.Range("M122").Value = -3797.0905
$ws.Range("N122").Value = -31797.31
$ws.Range("H136").Value = 13434.5
$ws.Range("I136").Value = 2799
$ws.Range("K136").Value = 8397
$ws.Range("M136").Value = -5847

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H123").Value = 59900
$ws.Range("J123").Value = 59900
$ws.Range("L123").Value = 59900
$ws.Range("N123").Value = -69700
$ws.Range("H134").Value = 2211.614
$ws.Range("I134").Value = 1833.8723
$ws.Range("J134").Value = 3987
$ws.Range("K134").Value = 5501.6169
$ws.Range("L134").Value = 11961
$ws.Range("M134").Value = -2966.6169
$ws.Range("N134").Value = -17031

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1672.2
$ws.Range("I31").Value = 850.119
$ws.Range("J31").Value = 4328.154
$ws.Range("K31").Value = 850.119
$ws.Range("L31").Value = 4328.154
$ws.Range("M31").Value = -555.119
$ws.Range("N31").Value = -4918.154
$ws.Range("H34").Value = 1672.2
$ws.Range("I34").Value = 850.119
$ws.Range("J34").Value = 4328.154
$ws.Range("K34").Value = 850.119
$ws.Range("L34").Value = 4328.154
$ws.Range("M34").Value = -648.119
$ws.Range("N34").Value = -4732.154
$ws.Range("H69").Value = 95995.4
$ws.Range("J69").Value = 114994.25
$ws.Range("L69").Value = 114994.25
$ws.Range("N69").Value = -116492.25
$ws.Range("H72").Value = 95995.4
$ws.Range("J72").Value = 114994.25
$ws.Range("L72").Value = 344982.75
$ws.Range("N72").Value = -352470.75
$ws.Range("H99").Value = 12413.866
$ws.Range("I99").Value = 20368.166
$ws.Range("J99").Value = 7111
$ws.Range("K99").Value = 20368.166
$ws.Range("L99").Value = 7111
$ws.Range("M99").Value = -18870.166
$ws.Range("N99").Value = -10107
$ws.Range("H126").Value = 12413.866
$ws.Range("I126").Value = 20368.166
$ws.Range("J126").Value = 7111
$ws.Range("K126").Value = 61104.49800000001
$ws.Range("L126").Value = 21333
$ws.Range("M126").Value = -58634.49800000001
$ws.Range("N126").Value = -26273
$ws.Range("H132").Value = 16670829
$ws.Range("I132").Value = 18522540
$ws.Range("K132").Value = 55567620
$ws.Range("M132").Value = -55565090

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 399.33334
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 399.33334
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 1198.00002
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -1422.00002
$ws.Range("H18").Value = 1159.5
$ws.Range("I18").Value = 1333
$ws.Range("J18").Value = 639
$ws.Range("K18").Value = 3999
$ws.Range("L18").Value = 1917
$ws.Range("M18").Value = -3830
$ws.Range("N18").Value = -2255
$ws.Range("H22").Value = 1160.1
$ws.Range("J22").Value = 900.25
$ws.Range("L22").Value = 2700.75
$ws.Range("N22").Value = -3038.75
$ws.Range("H27").Value = 1160.1
$ws.Range("J27").Value = 900.25
$ws.Range("L27").Value = 2700.75
$ws.Range("N27").Value = -2904.75
$ws.Range("H131").Value = 7939188
$ws.Range("J131").Value = 6063294.5
$ws.Range("L131").Value = 18189883.5
$ws.Range("N131").Value = -18199963.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 17032.666
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 17032.666
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 17032.666
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -18090.666
$ws.Range("H122").Value = 552274.25
$ws.Range("J122").Value = 1997.8334
$ws.Range("L122").Value = 5993.5002
$ws.Range("N122").Value = -10893.5002
$ws.Range("H123").Value = 53715.75
$ws.Range("J123").Value = 53715.75
$ws.Range("L123").Value = 53715.75
$ws.Range("N123").Value = -58615.75
$ws.Range("H126").Value = 4638.619
$ws.Range("I126").Value = 2746.25
$ws.Range("J126").Value = 7161.778
$ws.Range("K126").Value = 8238.75
$ws.Range("L126").Value = 21485.334
$ws.Range("M126").Value = -5768.75
$ws.Range("N126").Value = -26425.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 10000
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H40").Value = 5221
$ws.Range("I40").Value = 4181.4546
$ws.Range("K40").Value = 4181.4546
$ws.Range("M40").Value = -4045.4546
$ws.Range("H46").Value = 5626.515
$ws.Range("I46").Value = 3112.75
$ws.Range("K46").Value = 3112.75
$ws.Range("M46").Value = -2924.75
$ws.Range("H82").Value = 3908166.5
$ws.Range("I82").Value = 10418667
$ws.Range("J82").Value = 1866.4
$ws.Range("K82").Value = 10418667
$ws.Range("L82").Value = 1866.4
$ws.Range("M82").Value = -10418306
$ws.Range("N82").Value = -2588.4
$ws.Range("H85").Value = 3908166.5
$ws.Range("I85").Value = 10418667
$ws.Range("J85").Value = 1866.4
$ws.Range("K85").Value = 10418667
$ws.Range("L85").Value = 1866.4
$ws.Range("M85").Value = -10417419
$ws.Range("N85").Value = -4362.4
$ws.Range("H122").Value = 126988370
$ws.Range("I122").Value = 200003250
$ws.Range("K122").Value = 600009750
$ws.Range("M122").Value = -600007300
$ws.Range("H136").Value = 7257.2144
$ws.Range("I136").Value = 5433.6665
$ws.Range("K136").Value = 16300.9995
$ws.Range("M136").Value = -13750.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2093704.8
$ws.Range("I81").Value = 1393188.2
$ws.Range("K81").Value = 2786376.4
$ws.Range("M81").Value = -2785315.4
$ws.Range("H84").Value = 2093704.8
$ws.Range("I84").Value = 1393188.2
$ws.Range("K84").Value = 13931882
$ws.Range("M84").Value = -13926578
$ws.Range("H132").Value = 15628891
$ws.Range("I132").Value = 1672.3448
$ws.Range("K132").Value = 5017.0344
$ws.Range("M132").Value = -2487.0344
